$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new "banco de ditados" columns ---
$ws.Range("A1").Value = "id_ditado"
$ws.Range("B1").Value = "titulo"
$ws.Range("C1").Value = "traducaoTitulo"
$ws.Range("D1").Value = "descricao"
$ws.Range("E1").Value = "exemplo"
$ws.Range("F1").Value = "exemploTraducao"
$ws.Range("G1").Value = "dificuldade"

# --- Data row (row 2) ---
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "Put the horses before the cars"
$ws.Range("C2").Value = "Coloque os cavalos antes dos carros"
$ws.Range("D2").Value = "Nao se precipitar"
$ws.Range("E2").Value = "Put the horses before the cars"
$ws.Range("F2").Value = "Coloque os cavalos antes dos carros"
$ws.Range("G2").Value = 3

# New empty, underlined-style cell H2
$ws.Range("H2").Value = ""
$ws.Range("H2").Font.Underline = $true

# --- Sheet view / selection ---
$ws.Range("D4").Select

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
